$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '69.214.88'

# Row 3
Set-TextValue 'D3' '3.423.94'
Set-TextValue 'E3' '  +1.12%  '

# Row 4
Set-TextValue 'E4' '  +0.07%  '

# Row 5
Set-TextValue 'D5' '578.93'
Set-TextValue 'E5' '  -1.60%  '

# Row 6
Set-TextValue 'D6' '176.86'
Set-TextValue 'E6' '  -2.29%  '

# Row 7
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  +0.12%  '

# Row 8
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-TextValue 'D8' '3.415.76'
Set-TextValue 'E8' '  +1.09%  '

# Row 9
Set-TextValue 'D9' '0.590'
Set-TextValue 'E9' '  -0.96%  '

# Row 10
Set-TextValue 'D10' '0.196'
Set-TextValue 'E10' '  -0.06%  '

# Row 11
Set-TextValue 'D11' '0.583'
Set-TextValue 'E11' '  -0.93%  '

# Row 12
Set-TextValue 'D12' '48.56'
Set-TextValue 'E12' '  -0.93%  '

# Row 13
Set-TextValue 'D13' '0.0000278'
Set-TextValue 'E13' '  -1.64%  '

# Row 14
Set-TextValue 'D14' '694.64'
Set-TextValue 'E14' '  +1.16%  '

# Row 15
Set-TextValue 'D15' '3.956.77'
Set-TextValue 'E15' '  +0.41%  '

# Row 16
Set-TextValue 'D16' '8.60'
Set-TextValue 'E16' '  +0.05%  '

# Row 17
Set-TextValue 'D17' '69.313.86'
Set-TextValue 'E17' '  -0.17%  '

# Row 18
Set-TextValue 'D18' '3.428.88'
Set-TextValue 'E18' '  +1.41%  '

# Row 19
Set-TextValue 'E19' '  +0.72%  '

# Row 20
Set-TextValue 'D20' '17.68'
Set-TextValue 'E20' '  -0.40%  '

# Row 21
Set-TextValue 'D21' '11.35'
Set-TextValue 'E21' '  -0.47%  '

# Row 22
Set-TextValue 'D22' '0.896'
Set-TextValue 'E22' '  -0.97%  '

# Row 23
Set-TextValue 'E23' '  -0.04%  '

# Row 24
Set-TextValue 'D24' '16.91'
Set-TextValue 'E24' '  -1.00%  '

# Row 25
Set-TextValue 'D25' '100.51'
Set-TextValue 'E25' '  -4.00%  '

# Row 26
Set-TextValue 'D26' '3.88'
Set-TextValue 'E26' '  -2.05%  '

# Row 27
Set-TextValue 'D27' '2.66'
Set-TextValue 'E27' '  -2.52%  '

# Row 28
Set-TextValue 'D28' '9.56'
Set-TextValue 'E28' '  -0.65%  '

# Row 29
Set-TextValue 'D29' '33.42'
Set-TextValue 'E29' '  -3.68%  '

# Row 30
Set-TextValue 'D30' '8.71'
Set-TextValue 'E30' '  +0.04%  '

# Row 31
Set-TextValue 'D31' '6.92'
Set-TextValue 'E31' '  -1.71%  '

# Row 32
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D32' '567.77'
Set-TextValue 'E32' '  +2.18%  '

# Row 33
$ws.Range('B33').Value = 'dogwifhat'
$ws.Range('C33').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D33' '3.68'
Set-TextValue 'E33' '  +0.67%  '

# Row 34
Set-TextValue 'D34' '10.98'
Set-TextValue 'E34' '  -1.82%  '

# Row 35
Set-TextValue 'E35' '  -1.77%  '

# Row 36
Set-TextValue 'D36' '58.18'
Set-TextValue 'E36' '  -0.32%  '

# Row 37
Set-TextValue 'D37' '0.999'
Set-TextValue 'E37' '  -0.03%  '

# Row 38
Set-TextValue 'D38' '3.590.49'
Set-TextValue 'E38' '  -3.31%  '

# Row 39
Set-TextValue 'D39' '0.138'
Set-TextValue 'E39' '  -2.11%  '

# Row 40
Set-TextValue 'D40' '34.81'
Set-TextValue 'E40' '  -0.66%  '

# Row 41
Set-TextValue 'D41' '0.0₃0725'
Set-TextValue 'E41' '  +2.40%  '

# Row 42
Set-TextValue 'D42' '3.25'
Set-TextValue 'E42' '  -0.76%  '

# Row 43
Set-TextValue 'D43' '2.66'
Set-TextValue 'E43' '  -0.62%  '

# Row 44
Set-TextValue 'D44' '3.36'
Set-TextValue 'E44' '  +3.37%  '

# Row 45
Set-TextValue 'D45' '0.331'
Set-TextValue 'E45' '  -2.59%  '

# Row 46
Set-TextValue 'D46' '0.0418'
Set-TextValue 'E46' '  -0.31%  '

# Row 47
Set-TextValue 'D47' '1.46'
Set-TextValue 'E47' '  +3.26%  '

# Row 48
Set-TextValue 'D48' '2.64'
Set-TextValue 'E48' '  -0.59%  '

# Row 49
Set-TextValue 'E49' '  -1.59%  '

# Row 50
Set-TextValue 'E50' '  -0.13%  '

# Row 51
Set-TextValue 'D51' '131.31'
Set-TextValue 'E51' '  -1.12%  '
